$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values on the rows affected by the repulled data.
$ws.Range("F3").Value = 5
$ws.Range("F5").Value = -2
$ws.Range("F8").Value = -6
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = 3
$ws.Range("F17").Value = -2
$ws.Range("F19").Value = 3
$ws.Range("F23").Value = 0
$ws.Range("F31").Value = -1
$ws.Range("F35").Value = 4
$ws.Range("F38").Value = 3
$ws.Range("F39").Value = -1
$ws.Range("F44").Value = -1
$ws.Range("F47").Value = 1
$ws.Range("F52").Value = -1
$ws.Range("F54").Value = -1
$ws.Range("F58").Value = -1
$ws.Range("F62").Value = -1
$ws.Range("F63").Value = -3
$ws.Range("F64").Value = 3
$ws.Range("F65").Value = 0
$ws.Range("F67").Value = 5
